$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.349.24'
$ws.Range("E2").Value = '  -2.88%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.59'
$ws.Range("E5").Value = '  -3.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4240'
$ws.Range("E7").Value = '  -9.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3610'
$ws.Range("E8").Value = '  -2.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.29'
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07421'
$ws.Range("E10").Value = '  -3.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.113'
$ws.Range("E11").Value = '  -3.59%  '
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("E13").Value = '  -5.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.074'
$ws.Range("E14").Value = '  -4.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.176'
$ws.Range("E15").Value = '  -3.10%  '
$ws.Range("D16").Value = '1.730.36'
$ws.Range("E16").Value = '  -3.87%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001062'
$ws.Range("E17").Value = '  -3.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '87.64'
$ws.Range("E18").Value = '  +6.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06025'
$ws.Range("E19").Value = '  -10.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.82'
$ws.Range("E21").Value = '  -3.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.094'
$ws.Range("E22").Value = '  -5.28%  '
$ws.Range("E23").Value = '  -6.02%  '
$ws.Range("D24").Value = '27.386.02'
$ws.Range("E24").Value = '  -2.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.35'
$ws.Range("E25").Value = '  -4.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.383'
$ws.Range("E26").Value = '  -1.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.16'
$ws.Range("E27").Value = '  -3.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.372'
$ws.Range("E28").Value = '  -1.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '149.23'
$ws.Range("E29").Value = '  -2.43%  '
$ws.Range("D30").Value = '1.928.59'
$ws.Range("E30").Value = '  -3.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '126.40'
$ws.Range("E31").Value = '  -5.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.173'
$ws.Range("E32").Value = '  -7.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.669'
$ws.Range("E33").Value = '  -4.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09095'
$ws.Range("E34").Value = '  -5.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.625'
$ws.Range("E35").Value = '  -10.32%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '13.16'
$ws.Range("E36").Value = '  +7.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2141'
$ws.Range("E37").Value = '  -4.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.074'
$ws.Range("E38").Value = '  -3.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02253'
$ws.Range("E39").Value = '  -5.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06044'
$ws.Range("E40").Value = '  -5.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6377'
$ws.Range("E41").Value = '  -5.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.187'
$ws.Range("E42").Value = '  -4.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.967'
$ws.Range("E43").Value = '  -2.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.406'
$ws.Range("E45").Value = '  -7.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.60'
$ws.Range("E46").Value = '  -4.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.715'
$ws.Range("E47").Value = '  -3.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5821'
$ws.Range("E48").Value = '  -6.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.18'
$ws.Range("E49").Value = '  -3.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.957'
$ws.Range("E50").Value = '  -5.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06854'
